$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date only
$ws.Range("D2").Value = 44448

# Row 3: date, calidad, precios
$ws.Range("D3").Value = 44460
$ws.Range("L3").Value = "Especial"
$ws.Range("N3").Value = 31000
$ws.Range("O3").Value = 32000
$ws.Range("P3").Value = 31500
$ws.Range("S3").Value = 3150

# Row 4: date, precios
$ws.Range("D4").Value = 44460
$ws.Range("N4").Value = 30000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 30000
$ws.Range("S4").Value = 3000

# Row 5: date, calidad, precios
$ws.Range("D5").Value = 44446
$ws.Range("L5").Value = "Primera"
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("S5").Value = 2150

# Row 6: date, precios
$ws.Range("D6").Value = 44487
$ws.Range("N6").Value = 23000
$ws.Range("O6").Value = 24000
$ws.Range("P6").Value = 23500
$ws.Range("S6").Value = 2350

# Row 7: date, calidad, precios
$ws.Range("D7").Value = 44461
$ws.Range("L7").Value = "Especial"
$ws.Range("N7").Value = 31000
$ws.Range("O7").Value = 32000
$ws.Range("P7").Value = 31500
$ws.Range("S7").Value = 3150

# Row 8: date, volumen, precios
$ws.Range("D8").Value = 44461
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 30000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 30000
$ws.Range("S8").Value = 3000

# Row 9: date, calidad, precios
$ws.Range("D9").Value = 44452
$ws.Range("L9").Value = "Primera"
$ws.Range("N9").Value = 21000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 21500
$ws.Range("S9").Value = 2150

# Row 10: date, volumen, precios
$ws.Range("D10").Value = 44447
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 21000
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 21500
$ws.Range("S10").Value = 2150
